$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- FIX DSBS Import mengikuti daftar dari pusat ---
# Change the PENCACAH e-mail first so the old shared string slot
# (previously "pcl01@bpssumsel.com", used only by G2) is reused in place.
$ws.Range("G2").Value = "idris@bps.go.id"

# KD KEC / KD DESA / email now follow the central (pusat) list.
$ws.Range("B2").Value = "093"
$ws.Range("C2").Value = "001"
$ws.Range("D2").Value = "001"
$ws.Range("I2").Value = "susenas"

# NKS becomes a real number instead of text, formatted like G2 (General).
$ws.Range("G2").NumberFormat = "0"
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").NumberFormat = "0"
$ws.Range("E2").Value = 1601093001001

# --- styling ---
# Selection moves from I6 to I3.
$ws.Range("I3").Select()

# Narrow columns G and H.
$ws.Columns.Item(7).ColumnWidth = 16.83
$ws.Columns.Item(8).ColumnWidth = 8.5
